$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# QR check-in updates: the reservation name cells for rows 4 and 5 are
# reassigned (date/status stay the same) -- row 4's reservation is now
# "monji" and row 5's reservation is now "ahmed".
$ws.Range("A4").Value = "monji"
$ws.Range("A5").Value = "ahmed"
